# Add the 17-Feb-2020 timesheet entries to the bottom of the sheet.
#
# The sheet is laid out as repeating "days": a highlighted blank separator
# row (A/C centered, B left-aligned, orange/theme fill) followed by one
# row per task with columns Timestamp | Task | Location. Rows whose Task
# text wraps to two lines use a taller (30pt) row and a wrap-aligned style
# for column B.
#
# NB: this interpreter's named-parameter binding for user functions is
# unreliable, so the per-row logic below is written out inline instead of
# being wrapped in a helper function.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Blank separator row (same look as the row above every other day) ---
$ws.Range("A132:C132").Copy()
$ws.Range("A140:C140").PasteSpecial(-4122)

# --- 2) Data rows for Feb 17 ------------------------------------------------
# Default (single line) row style template: A139/B139/C139
# Wrapped (two line) row style template:     A133/B133/C133

# Row 141
$ws.Range("A139:C139").Copy()
$ws.Range("A141:C141").PasteSpecial(-4122)
$ws.Range("A141").Value = "Feb 17 10:00 to 11:00"
$ws.Range("B141").Value = "Modified code of hourly data buckets, added combined data check"
$ws.Range("C141").Value = "Infimetrics"

# Row 142
$ws.Range("A139:C139").Copy()
$ws.Range("A142:C142").PasteSpecial(-4122)
$ws.Range("A142").Value = "Feb 17 11:00 to 12:00"
$ws.Range("B142").Value = "Documenting and commenting HourWiseData class"
$ws.Range("C142").Value = "Infimetrics"

# Row 143 (wraps to two lines)
$ws.Range("A133:C133").Copy()
$ws.Range("A143:C143").PasteSpecial(-4122)
$ws.Range("A143").Value = "Feb 17 12:00 to 13:00"
$ws.Range("B143").Value = "Documented and commented HourWiseData class. Documenting and commenting" + [char]10 + "HourlyMachineData class"
$ws.Range("C143").Value = "Infimetrics"
$ws.Rows.Item(143).RowHeight = 30

# Row 144
$ws.Range("A139:C139").Copy()
$ws.Range("A144:C144").PasteSpecial(-4122)
$ws.Range("A144").Value = "Feb 17 13:00 to 13:30"
$ws.Range("B144").Value = "Documented and commented HourlyMachineData class."
$ws.Range("C144").Value = "Infimetrics"

# Row 145
$ws.Range("A139:C139").Copy()
$ws.Range("A145:C145").PasteSpecial(-4122)
$ws.Range("A145").Value = "Feb 17 13:30 to 14:00"
$ws.Range("B145").Value = "Lunch"
$ws.Range("C145").Value = "Infimetrics"

# Row 146
$ws.Range("A139:C139").Copy()
$ws.Range("A146:C146").PasteSpecial(-4122)
$ws.Range("A146").Value = "Feb 17 14:00 to 15:00"
$ws.Range("B146").Value = "Data Cleaning done. Working on data transformation."
$ws.Range("C146").Value = "Infimetrics"

# Row 147
$ws.Range("A139:C139").Copy()
$ws.Range("A147:C147").PasteSpecial(-4122)
$ws.Range("A147").Value = "Feb 17 15:00 to 16:00"
$ws.Range("B147").Value = "Getting errors and complication in data transformation. Working on fixing issues"
$ws.Range("C147").Value = "Infimetrics"

# Row 148 (wraps to two lines)
$ws.Range("A133:C133").Copy()
$ws.Range("A148:C148").PasteSpecial(-4122)
$ws.Range("A148").Value = "Feb 17 16:00 to 17:00"
$ws.Range("B148").Value = "Not considering data transformation due to issues in loading saved combined data." + [char]10 + "Removed code of saving combined data."
$ws.Range("C148").Value = "Infimetrics"
$ws.Rows.Item(148).RowHeight = 30

# Row 149 (wraps to two lines)
$ws.Range("A133:C133").Copy()
$ws.Range("A149:C149").PasteSpecial(-4122)
$ws.Range("A149").Value = "Feb 17 17:00 to 18:00"
$ws.Range("B149").Value = "Worked on feature engineering. Created features as max occurred sublocation, total" + [char]10 + "alarm duration and total automation duration."
$ws.Range("C149").Value = "Infimetrics"
$ws.Rows.Item(149).RowHeight = 30

# Row 150
$ws.Range("A139:C139").Copy()
$ws.Range("A150:C150").PasteSpecial(-4122)
$ws.Range("A150").Value = "Feb 17 18:00 to 19:00"
$ws.Range("B150").Value = "Getting issues and complications in feature engineering progress bars."
$ws.Range("C150").Value = "Infimetrics"

# --- 3) Scroll / selection, matching where the author ended up -------------
$ws.Application.Goto($ws.Range("A134"), $true)
$ws.Application.Goto($ws.Range("D150"))
